$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Principal Component Analysis")

# --- Restructure the existing "Python" table: shift it from A:C to B:D,
# and push everything down one row so a new header row (row 1) can hold
# the "Python" / "Orange" / "Data Polish" titles (mirrors the
# "Numeric To Categorical" sheet's layout).
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# Content column (now column D) gets the wrap-text variant style used for
# the "Content" column elsewhere in the workbook.
$ws.Range("D2:D6").WrapText = $true

# Row heights for the (now taller, wrapped) content rows.
$ws.Rows.Item(2).RowHeight = 18
$ws.Rows.Item(3).RowHeight = 72
$ws.Rows.Item(4).RowHeight = 108
$ws.Rows.Item(5).RowHeight = 126
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 17

# --- New header row: tool names + overall section title.
$ws.Range("C1").Value = "Python"
$ws.Range("G1").Value = "Orange"
$ws.Range("K1").Value = "Data Polish"

foreach ($addr in @("C1", "D1", "G1", "H1", "K1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Font.Size = 14
}
foreach ($addr in @("D1", "H1")) {
    $ws.Range($addr).WrapText = $true
}
$ws.Rows.Item(1).RowHeight = 19

# --- New "Orange" table (PCA via Orange widgets), columns F:H.
$ws.Range("F2").Value = "Action"
$ws.Range("G2").Value = "Time"
$ws.Range("H2").Value = "Content"

$ws.Range("F3").Value = "Load Data"
$ws.Range("G3").Value = "2 min"
$ws.Range("H3").Value = "Use 'File' widget to load the dataset"

$ws.Range("F4").Value = "Apply PCA"
$ws.Range("G4").Value = "3 min"
$ws.Range("H4").Value = "Connect 'PCA' widget and set parameters"

$ws.Range("F5").Value = "Visualize Results"
$ws.Range("G5").Value = "2 min"
$ws.Range("H5").Value = "Connect 'Scatter Plot' widget to visualize PCA results"

$ws.Range("F6").Value = "Overall"
$ws.Range("G6").Value = "7 min"

# Header / footer rows (bold, Mac system font).
foreach ($addr in @("F2", "G2", "H2", "F6", "G6")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Font.Size = 13
    $c.Font.Name = ".AppleSystemUIFont"
}
$ws.Range("H2").WrapText = $true

# Data rows (regular weight, Mac system font, wrapped content column).
foreach ($addr in @("F3", "G3", "F4", "G4", "F5", "G5")) {
    $c = $ws.Range($addr)
    $c.Font.Size = 13
    $c.Font.Name = ".AppleSystemUIFont"
}
foreach ($addr in @("H3", "H4", "H5")) {
    $c = $ws.Range($addr)
    $c.Font.Size = 13
    $c.Font.Name = ".AppleSystemUIFont"
    $c.WrapText = $true
}

# Trailing empty cell under the Orange "Content" column, styled like the
# matching cell on the "Numeric To Categorical" sheet.
$h6 = $ws.Range("H6")
$h6.Font.Name = "Helvetica"
$h6.Font.Size = 12
$h6.WrapText = $true

# --- Column widths for the two "Content" columns (D and H).
$ws.Columns.Item(4).ColumnWidth = 10.83203125
$ws.Columns.Item(8).ColumnWidth = 10.83203125

# --- Selection bookkeeping: this sheet becomes the active tab, while the
# previously active "Numeric To Categorical" sheet selects its header row.
$ws5 = $wb.Worksheets.Item("Numeric To Categorical")
$ws5.Rows.Item(1).Select()

$ws.Activate()
$ws.Range("H4").Select()
